$p = $ppt.ActivePresentation

# Remove slides 2, 3, and 4 (keeping only the first slide), deleting from the
# end so indices remain valid as we go.
for ($i = $p.Slides.Count; $i -ge 2; $i--) {
    $p.Slides.Item($i).Delete()
}
